$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.250039333333334
$ws.Range("H2").Value = 6.750118000000001
$ws.Range("I2").Value = 0.05266340474459803
$ws.Range("J2").Value = 0.05266340474459803
$ws.Range("M2").Value = 0.9956583333333334
$ws.Range("N2").Value = 2.986975
$ws.Range("O2").Value = 0.1535710477437721
$ws.Range("P2").Value = 0.153571047743772
$ws.Range("Q2").Value = 2.240270412561112
$ws.Range("R2").Value = 20.16243371305
$ws.Range("S2").Value = 0.008087574244382256
$ws.Range("T2").Value = 0.008087574244382256
$ws.Range("G3").Value = 2.250039333333334
$ws.Range("H3").Value = 6.750118000000001
$ws.Range("I3").Value = 0.05266340474459803
$ws.Range("J3").Value = 0.05266340474459803
$ws.Range("O3").Value = 0.2875834480798523
$ws.Range("P3").Value = 0.2875834480798522
$ws.Range("Q3").Value = 4.195222337419556
$ws.Range("R3").Value = 37.757001036776
$ws.Range("S3").Value = 0.01514512352407636
$ws.Range("T3").Value = 0.01514512352407635
$ws.Range("G4").Value = 2.250039333333334
$ws.Range("H4").Value = 6.750118000000001
$ws.Range("I4").Value = 0.05266340474459803
$ws.Range("J4").Value = 0.05266340474459803
$ws.Range("M4").Value = 2.910118
$ws.Range("N4").Value = 8.730354
$ws.Range("O4").Value = 0.4488586650219809
$ws.Range("P4").Value = 0.4488586650219808
$ws.Range("Q4").Value = 6.547879964641334
$ws.Range("R4").Value = 58.930919681772
$ws.Range("S4").Value = 0.02363842554917253
$ws.Range("T4").Value = 0.02363842554917252
$ws.Range("G5").Value = 2.250039333333334
$ws.Range("H5").Value = 6.750118000000001
$ws.Range("I5").Value = 0.05266340474459803
$ws.Range("J5").Value = 0.05266340474459803
$ws.Range("M5").Value = 0.7130856666666667
$ws.Range("N5").Value = 2.139257
$ws.Range("O5").Value = 0.1099868391543949
$ws.Range("P5").Value = 0.1099868391543948
$ws.Range("Q5").Value = 1.604470798036223
$ws.Range("R5").Value = 14.440237182326
$ws.Range("S5").Value = 0.0057922814269669
$ws.Range("T5").Value = 0.005792281426966897
$ws.Range("I6").Value = 0.4446323259834918
$ws.Range("J6").Value = 0.4446323259834918
$ws.Range("M6").Value = 0.9956583333333334
$ws.Range("N6").Value = 2.986975
$ws.Range("O6").Value = 0.1535710477437721
$ws.Range("P6").Value = 0.153571047743772
$ws.Range("Q6").Value = 18.91439889235834
$ws.Range("R6").Value = 170.229590031225
$ws.Range("S6").Value = 0.06828265216203525
$ws.Range("T6").Value = 0.06828265216203523
$ws.Range("I7").Value = 0.4446323259834918
$ws.Range("J7").Value = 0.4446323259834918
$ws.Range("O7").Value = 0.2875834480798523
$ws.Range("P7").Value = 0.2875834480798522
$ws.Range("S7").Value = 0.1278688974340975
$ws.Range("T7").Value = 0.1278688974340974
$ws.Range("I8").Value = 0.4446323259834918
$ws.Range("J8").Value = 0.4446323259834918
$ws.Range("M8").Value = 2.910118
$ws.Range("N8").Value = 8.730354
$ws.Range("O8").Value = 0.4488586650219809
$ws.Range("P8").Value = 0.4488586650219808
$ws.Range("Q8").Value = 55.28315370148601
$ws.Range("R8").Value = 497.5483833133741
$ws.Range("S8").Value = 0.1995770722665683
$ws.Range("T8").Value = 0.1995770722665683
$ws.Range("I9").Value = 0.4446323259834918
$ws.Range("J9").Value = 0.4446323259834918
$ws.Range("M9").Value = 0.7130856666666667
$ws.Range("N9").Value = 2.139257
$ws.Range("O9").Value = 0.1099868391543949
$ws.Range("P9").Value = 0.1099868391543948
$ws.Range("Q9").Value = 13.54640070012967
$ws.Range("R9").Value = 121.917606301167
$ws.Range("S9").Value = 0.04890370412079078
$ws.Range("T9").Value = 0.04890370412079077
$ws.Range("G10").Value = 21.38027566666667
$ws.Range("H10").Value = 64.140827
$ws.Range("I10").Value = 0.5004170790724312
$ws.Range("J10").Value = 0.5004170790724313
$ws.Range("M10").Value = 0.9956583333333334
$ws.Range("N10").Value = 2.986975
$ws.Range("O10").Value = 0.1535710477437721
$ws.Range("P10").Value = 0.153571047743772
$ws.Range("Q10").Value = 21.28744963648056
$ws.Range("R10").Value = 191.587046728325
$ws.Range("S10").Value = 0.07684957514203128
$ws.Range("T10").Value = 0.07684957514203129
$ws.Range("G11").Value = 21.38027566666667
$ws.Range("H11").Value = 64.140827
$ws.Range("I11").Value = 0.5004170790724312
$ws.Range("J11").Value = 0.5004170790724313
$ws.Range("O11").Value = 0.2875834480798523
$ws.Range("P11").Value = 0.2875834480798522
$ws.Range("Q11").Value = 39.86375203677377
$ws.Range("R11").Value = 358.773768330964
$ws.Range("S11").Value = 0.1439116690776979
$ws.Range("T11").Value = 0.1439116690776978
$ws.Range("G12").Value = 21.38027566666667
$ws.Range("H12").Value = 64.140827
$ws.Range("I12").Value = 0.5004170790724312
$ws.Range("J12").Value = 0.5004170790724313
$ws.Range("M12").Value = 2.910118
$ws.Range("N12").Value = 8.730354
$ws.Range("O12").Value = 0.4488586650219809
$ws.Range("P12").Value = 0.4488586650219808
$ws.Range("Q12").Value = 62.21912506252867
$ws.Range("R12").Value = 559.972125562758
$ws.Range("S12").Value = 0.2246165420666505
$ws.Range("T12").Value = 0.2246165420666505
$ws.Range("G13").Value = 21.38027566666667
$ws.Range("H13").Value = 64.140827
$ws.Range("I13").Value = 0.5004170790724312
$ws.Range("J13").Value = 0.5004170790724313
$ws.Range("M13").Value = 0.7130856666666667
$ws.Range("N13").Value = 2.139257
$ws.Range("O13").Value = 0.1099868391543949
$ws.Range("P13").Value = 0.1099868391543948
$ws.Range("Q13").Value = 15.24596812728211
$ws.Range("R13").Value = 137.213713145539
$ws.Range("S13").Value = 0.05503929278605159
$ws.Range("T13").Value = 0.05503929278605158
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.09771999999999999
$ws.Range("H14").Value = 0.29316
$ws.Range("I14").Value = 0.002287190199478936
$ws.Range("J14").Value = 0.002287190199478936
$ws.Range("M14").Value = 0.9956583333333334
$ws.Range("N14").Value = 2.986975
$ws.Range("O14").Value = 0.1535710477437721
$ws.Range("P14").Value = 0.153571047743772
$ws.Range("Q14").Value = 0.09729573233333333
$ws.Range("R14").Value = 0.875661591
$ws.Range("S14").Value = 0.0003512461953232672
$ws.Range("T14").Value = 0.0003512461953232672
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.09771999999999999
$ws.Range("H15").Value = 0.29316
$ws.Range("I15").Value = 0.002287190199478936
$ws.Range("J15").Value = 0.002287190199478936
$ws.Range("O15").Value = 0.2875834480798523
$ws.Range("P15").Value = 0.2875834480798522
$ws.Range("Q15").Value = 0.1821999823466667
$ws.Range("R15").Value = 1.63979984112
$ws.Range("S15").Value = 0.0006577580439805975
$ws.Range("T15").Value = 0.0006577580439805974
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.09771999999999999
$ws.Range("H16").Value = 0.29316
$ws.Range("I16").Value = 0.002287190199478936
$ws.Range("J16").Value = 0.002287190199478936
$ws.Range("M16").Value = 2.910118
$ws.Range("N16").Value = 8.730354
$ws.Range("O16").Value = 0.4488586650219809
$ws.Range("P16").Value = 0.4488586650219808
$ws.Range("Q16").Value = 0.28437673096
$ws.Range("R16").Value = 2.55939057864
$ws.Range("S16").Value = 0.001026625139589473
$ws.Range("T16").Value = 0.001026625139589473
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.09771999999999999
$ws.Range("H17").Value = 0.29316
$ws.Range("I17").Value = 0.002287190199478936
$ws.Range("J17").Value = 0.002287190199478936
$ws.Range("M17").Value = 0.7130856666666667
$ws.Range("N17").Value = 2.139257
$ws.Range("O17").Value = 0.1099868391543949
$ws.Range("P17").Value = 0.1099868391543948
$ws.Range("Q17").Value = 0.06968273134666667
$ws.Range("R17").Value = 0.62714458212
$ws.Range("S17").Value = 0.000251560820585598
$ws.Range("T17").Value = 0.000251560820585598
